$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 53-59: O column tied_teams - move 'Ireland' from front to back
foreach ($r in 53..59) {
    $ws.Range("O$r").Value = "['Colombia', 'Costa Rica', 'Argentina', 'Ireland']"
}

# Rows 63-73: O column tied_teams - move 'Scotland' from position 2 to back
foreach ($r in 63..73) {
    $ws.Range("O$r").Value = "['Colombia', 'Austria', 'Argentina', 'Scotland']"
}

# Row 78: O column tied_teams - swap order
$ws.Range("O78").Value = "['South Korea', 'Netherlands']"

# Row 108: J, M, P, Q updates (Bulgaria -> Argentina tie break)
$ws.Range("J108").Value = "['Argentina', 6, 4, 6]"
$ws.Range("M108").Value = "['Argentina', 'Belgium', 'United States', 'Italy']"
$ws.Range("P108").Value = 1
$ws.Range("Q108").Value = 12

# Row 109: P updates
$ws.Range("P109").Value = 0
